# Update sprint back log for cooper
# Adds two new sprint-backlog rows (10 & 11) for Cooper Dahlberg on the
# active sheet ("Sheet2"), matching the formatting of the existing data
# rows, and moves the sheet's active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting -----------------------------------------------------
# Row 9 (A:F) already carries the "data row" styling we want for the two
# new rows, so copy it down. Column G needs the style used by rows
# 5-7 (border only, no forced left-alignment) rather than row 9's G style.
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)
$ws.Range("A9:F9").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("G11").PasteSpecial(-4122)

# --- Row 10: "lose lives" story --------------------------------------
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "Cooper Dahlberg"
$ws.Range("D10").Value = "F"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = "As a player I want to be able to lose lives so that the game has risks in it."

# --- Row 11: "gameover" story -----------------------------------------
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Cooper Dahlberg"
$ws.Range("D11").Value = "F"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = "As a player I want to be able to get a gameover so that I can lose the game."

# --- Selection ----------------------------------------------------------
[void]$ws.Range("F19").Select()
